# Update countries & provincias Spain
# - Re-sort a few country rows whose case counts changed order
#   (Rumania/Catar, Nepal/Honduras, Gibraltar/San Martin, Santa Lucia/Timor Oriental)
# - Refresh statistics (B:H) for the affected + several other rows
# - Refresh the "Datos actualizados..." timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 29 de Septiembre de 2020 a las 13:09'
$ws.Range("B4").Value = 7361889
$ws.Range("C4").Value = 278
$ws.Range("D4").Value = 4610639
$ws.Range("E4").Value = 2541435
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 209815
$ws.Range("B25").Value = 288745
$ws.Range("C25").Value = 127
$ws.Range("E25").Value = 26800
$ws.Range("A32").Value = 'Rumania'
$ws.Range("B32").Value = 125414
$ws.Range("C32").Value = 1470
$ws.Range("D32").Value = 100636
$ws.Range("E32").Value = 19986
$ws.Range("G32").Value = 44
$ws.Range("H32").Value = 4792
$ws.Range("A33").Value = 'Catar'
$ws.Range("B33").Value = 125311
$ws.Range("D33").Value = 122209
$ws.Range("E33").Value = 2888
$ws.Range("H33").Value = 214
$ws.Range("B42").Value = 98585
$ws.Range("C42").Value = 528
$ws.Range("D42").Value = 88528
$ws.Range("E42").Value = 9122
$ws.Range("G42").Value = 11
$ws.Range("H42").Value = 935
$ws.Range("B49").Value = 78260
$ws.Range("C49").Value = 314
$ws.Range("D49").Value = 74320
$ws.Range("E49").Value = 3112
$ws.Range("G49").Value = 6
$ws.Range("H49").Value = 828
$ws.Range("A50").Value = 'Nepal'
$ws.Range("B50").Value = 76258
$ws.Range("C50").Value = 1513
$ws.Range("D50").Value = 55371
$ws.Range("E50").Value = 20396
$ws.Range("G50").Value = 10
$ws.Range("H50").Value = 491
$ws.Range("A51").Value = 'Honduras'
$ws.Range("B51").Value = 75537
$ws.Range("C51").Value = 428
$ws.Range("D51").Value = 26957
$ws.Range("E51").Value = 46279
$ws.Range("G51").Value = 12
$ws.Range("H51").Value = 2301
$ws.Range("B61").Value = 52871
$ws.Range("C61").Value = 225
$ws.Range("E61").Value = 8102
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 2069
$ws.Range("B67").Value = 44041
$ws.Range("C67").Value = 609
$ws.Range("D67").Value = 34916
$ws.Range("E67").Value = 8329
$ws.Range("G67").Value = 6
$ws.Range("H67").Value = 796
$ws.Range("B91").Value = 14945
$ws.Range("C91").Value = 26
$ws.Range("D91").Value = 12335
$ws.Range("E91").Value = 2300
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 310
$ws.Range("B99").Value = 10631
$ws.Range("C99").Value = 7
$ws.Range("D99").Value = 10129
$ws.Range("E99").Value = 230
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 272
$ws.Range("B146").Value = 3035
$ws.Range("C146").Value = 29
$ws.Range("D146").Value = 2484
$ws.Range("E146").Value = 517
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = 34
$ws.Range("A180").Value = 'Gibraltar'
$ws.Range("B180").Value = 391
$ws.Range("C180").Value = 9
$ws.Range("D180").Value = 342
$ws.Range("E180").Value = 49
$ws.Range("H180").Value = 0
$ws.Range("A181").Value = 'San Martin (Parte Francesa)'
$ws.Range("B181").Value = 383
$ws.Range("D181").Value = 273
$ws.Range("E181").Value = 102
$ws.Range("H181").Value = 8
$ws.Range("A207").Value = 'Santa Lucia'
$ws.Range("A208").Value = 'Timor Oriental'
